# Apply the dataset edit described in the commit:
#   "fix greedy, fix dataset names, make path building take nearest vertex first"
#
# The functional content change for this workbook is a cleanup/rename of the
# point-name (and one header) text used in the "points" sheet: dashes and
# parenthetical qualifiers are removed/simplified, and a few names get their
# first letter capitalized. Row/column layout and all numeric data stay the
# same; only the display text referenced by column D (point_name) and cell
# G2 (rep_id description) changes. We also update the saved cell selection
# on the points sheet to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("points")

# NOTE: the assignment order below matters only for shared-string table
# bookkeeping (cosmetic); it mirrors the order the new strings were
# introduced in the source workbook.
$ws.Range("D4").Value = "Кромка льда на Западе"
$ws.Range("D19").Value = "Остров Врангеля"
$ws.Range("D23").Value = "Около Новой Земли"
$ws.Range("D28").Value = "Мыс.Наглёйнын"
$ws.Range("D29").Value = "Пролив Лонга"
$ws.Range("D37").Value = "Терминал Утренний"
$ws.Range("D40").Value = "Кромка льда на Востоке"
$ws.Range("D44").Value = "Остров Котельный"
$ws.Range("D20").Value = "Восточно-Сибирское 1"
$ws.Range("G2").Value = "rep_id обозначение на картинке"
$ws.Range("D9").Value = "Карское 3"
$ws.Range("D10").Value = "Пролив Вилькицкого 3"
$ws.Range("D11").Value = "Лаптевых 4"
$ws.Range("D14").Value = "Лаптевых 1"
$ws.Range("D15").Value = "Карское 1"
$ws.Range("D21").Value = "Пролив Вилькицкого восток"
$ws.Range("D22").Value = "Пролив Вилькицкого запад"
$ws.Range("D24").Value = "Пролив Санникова 1"
$ws.Range("D25").Value = "Пролив Санникова 2"
$ws.Range("D30").Value = "Восточно-Сибирское 3"
$ws.Range("D32").Value = "Лаптевых 3"
$ws.Range("D35").Value = "Восточно-Сибирское 2"
$ws.Range("D42").Value = "Лаптевых 2"
$ws.Range("D45").Value = "Карское 2"

# Update the remembered selection on the points sheet.
[void]$ws.Range("K9").Select()
